$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Domande")

$ws.Range("T2").Value = 1
$ws.Range("T5").Value = 2
$ws.Range("T33").Value = 3
$ws.Range("T61").Value = 4
$ws.Range("T74").Value = 5
$ws.Range("U2").Value = 1
$ws.Range("U5").Value = 2
$ws.Range("U14").Value = 3
$ws.Range("U17").Value = 4
$ws.Range("U19").Value = 5
$ws.Range("U26").Value = 6
$ws.Range("U29").Value = 7
$ws.Range("U33").Value = 8
$ws.Range("U36").Value = 9
$ws.Range("U61").Value = 10
$ws.Range("U65").Value = 11
$ws.Range("U68").Value = 12
$ws.Range("U70").Value = 13
$ws.Range("U72").Value = 14
$ws.Range("U74").Value = 15
$ws.Range("U76").Value = 16
$ws.Range("U78").Value = 17
$ws.Range("U80").Value = 18
$ws.Range("U82").Value = 19
$ws.Range("U84").Value = 20
$ws.Range("U86").Value = 21
$ws.Range("U90").Value = 22
$ws.Range("U92").Value = 23
$ws.Range("V2").Value = 1
$ws.Range("V3").Value = 2
$ws.Range("V4").Value = 3
$ws.Range("V5").Value = 4
$ws.Range("V13").Value = 5
$ws.Range("V14").Value = 6
$ws.Range("V17").Value = 7
$ws.Range("V19").Value = 8
$ws.Range("V20").Value = 9
$ws.Range("V21").Value = 10
$ws.Range("V22").Value = 11
$ws.Range("V23").Value = 12
$ws.Range("V24").Value = 13
$ws.Range("V25").Value = 14
$ws.Range("V26").Value = 15
$ws.Range("V27").Value = 16
$ws.Range("V28").Value = 17
$ws.Range("V29").Value = 18
$ws.Range("V33").Value = 19
$ws.Range("V34").Value = 20
$ws.Range("V35").Value = 21
$ws.Range("V36").Value = 22
$ws.Range("V44").Value = 23
$ws.Range("V45").Value = 24
$ws.Range("V53").Value = 25
$ws.Range("V61").Value = 26
$ws.Range("V62").Value = 27
$ws.Range("V63").Value = 28
$ws.Range("V64").Value = 29
$ws.Range("V65").Value = 30
$ws.Range("V66").Value = 31
$ws.Range("V67").Value = 32
$ws.Range("V68").Value = 33
$ws.Range("V70").Value = 34
$ws.Range("V72").Value = 35
$ws.Range("V74").Value = 36
$ws.Range("V76").Value = 37
$ws.Range("V78").Value = 38
$ws.Range("V80").Value = 39
$ws.Range("V82").Value = 40
$ws.Range("V84").Value = 41
$ws.Range("V86").Value = 42
$ws.Range("V90").Value = 43
$ws.Range("V92").Value = 44
$ws.Range("W2").Value = 1
$ws.Range("W3").Value = 2
$ws.Range("W4").Value = 3
$ws.Range("W5").Value = 4
$ws.Range("W6").Value = 5
$ws.Range("W7").Value = 6
$ws.Range("W8").Value = 7
$ws.Range("W9").Value = 8
$ws.Range("W10").Value = 9
$ws.Range("W11").Value = 10
$ws.Range("W12").Value = 11
$ws.Range("W13").Value = 12
$ws.Range("W14").Value = 13
$ws.Range("W15").Value = 14
$ws.Range("W16").Value = 15
$ws.Range("W17").Value = 16
$ws.Range("W18").Value = 17
$ws.Range("W19").Value = 18
$ws.Range("W20").Value = 19
$ws.Range("W21").Value = 20
$ws.Range("W22").Value = 21
$ws.Range("W23").Value = 22
$ws.Range("W24").Value = 23
$ws.Range("W25").Value = 24
$ws.Range("W26").Value = 25
$ws.Range("W27").Value = 26
$ws.Range("W28").Value = 27
$ws.Range("W29").Value = 28
$ws.Range("W30").Value = 29
$ws.Range("W31").Value = 30
$ws.Range("W32").Value = 31
$ws.Range("W33").Value = 32
$ws.Range("W34").Value = 33
$ws.Range("W35").Value = 34
$ws.Range("W36").Value = 35
$ws.Range("W37").Value = 36
$ws.Range("W38").Value = 37
$ws.Range("W39").Value = 38
$ws.Range("W40").Value = 39
$ws.Range("W41").Value = 40
$ws.Range("W42").Value = 41
$ws.Range("W43").Value = 42
$ws.Range("W44").Value = 43
$ws.Range("W45").Value = 44
$ws.Range("W46").Value = 45
$ws.Range("W47").Value = 46
$ws.Range("W48").Value = 47
$ws.Range("W49").Value = 48
$ws.Range("W50").Value = 49
$ws.Range("W51").Value = 50
$ws.Range("W52").Value = 51
$ws.Range("W53").Value = 52
$ws.Range("W54").Value = 53
$ws.Range("W55").Value = 54
$ws.Range("W56").Value = 55
$ws.Range("W57").Value = 56
$ws.Range("W58").Value = 57
$ws.Range("W59").Value = 58
$ws.Range("W60").Value = 59
$ws.Range("W61").Value = 60
$ws.Range("W62").Value = 61
$ws.Range("W63").Value = 62
$ws.Range("W64").Value = 63
$ws.Range("W65").Value = 64
$ws.Range("W66").Value = 65
$ws.Range("W67").Value = 66
$ws.Range("W68").Value = 67
$ws.Range("W69").Value = 68
$ws.Range("W70").Value = 69
$ws.Range("W71").Value = 70
$ws.Range("W72").Value = 71
$ws.Range("W73").Value = 72
$ws.Range("W74").Value = 73
$ws.Range("W75").Value = 74
$ws.Range("W76").Value = 75
$ws.Range("W77").Value = 76
$ws.Range("W78").Value = 77
$ws.Range("W79").Value = 78
$ws.Range("W80").Value = 79
$ws.Range("W81").Value = 80
$ws.Range("W82").Value = 81
$ws.Range("W83").Value = 82
$ws.Range("W84").Value = 83
$ws.Range("W85").Value = 84
$ws.Range("W86").Value = 85
$ws.Range("W87").Value = 86
$ws.Range("W88").Value = 87
$ws.Range("W89").Value = 88
$ws.Range("W90").Value = 89
$ws.Range("W91").Value = 90
$ws.Range("W92").Value = 91
$ws.Range("W93").Value = 92
